$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new product row (row 3) mirroring the existing Prod1 row structure
$ws.Range("A3").Value = "Prod2"
$ws.Range("B3").Value = "ProdDesc2"
$ws.Range("C3").Value = "www.prod2.com"
$ws.Range("D3").Value = 500

# Apply same style as C2 (Hyperlink style) to C3
$ws.Range("C3").Style = "Hyperlink"

# Add hyperlink for the new product url cell
$ws.Hyperlinks.Add($ws.Range("C3"), "http://www.prod2.com/", "", "", "www.prod2.com")

# Update selection to match target state
$ws.Range("D3").Select()
